$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 08:27"

# --- Reorder country names that moved position ---
# Rumania / Bielorrusia swap (rows 42-43)
$ws.Range("A42").Value = "Bielorrusia"
$ws.Range("A43").Value = "Rumania"

# Belice inserted before Guadalupe, shifting Guadalupe/Burundi/Comoras down one row (rows 169-172)
$ws.Range("A169").Value = "Belice"
$ws.Range("A170").Value = "Guadalupe"
$ws.Range("A171").Value = "Burundi"
$ws.Range("A172").Value = "Comoras"

# Montserrat / Islas Malvinas swap (rows 213-214)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Update case numbers (columns B:H) ---

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# Pakistan (row 18)
Set-Row 18 288717 670 266301 16248 0 6 6168

# Israel (row 33)
Set-Row 33 92343 110 68395 23274 0 0 674

# Bielorrusia (row 42, formerly Rumania's row)
Set-Row 42 69424 0 66747 2070 0 0 607

# Rumania (row 43, formerly Bielorrusia's row)
Set-Row 43 69374 0 32334 34086 0 0 2954

# Kirguistan (row 55)
Set-Row 55 41856 211 34276 6085 0 2 1495

# El Salvador (row 73)
Set-Row 73 22619 0 10647 11360 0 9 612

# Belice (row 169, new data)
Set-Row 169 452 64 35 414 0 0 3

# Guadalupe (row 170)
Set-Row 170 446 0 289 143 0 0 14

# Burundi (row 171)
Set-Row 171 412 0 315 96 0 0 1

# Comoras (row 172)
Set-Row 172 403 0 379 17 0 0 7

# Islas Malvinas (row 213)
Set-Row 213 13 0 13 0 0 0 0

# Montserrat (row 214)
Set-Row 214 13 0 12 0 0 0 1
